$d = $word.ActiveDocument

# Green color used elsewhere in the document for "working" components (00B050 -> OLE BGR)
$green = 5287936

# 1) "De buzzer: op GPIO1 of 2" paragraph turns from red to green (buzzer now works).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "De buzzer*") {
        $p.Range.Font.Color = $green
        break
    }
}

# 2) "Led strip ... RGB_LED (JESPER rood groen blauw)" paragraph:
#    - remove the trailing " (JESPER rood groen blauw)" remark
#    - RGB_LED itself is now working -> green
$rng = $d.Content
$rng.Find.Execute(" (JESPER rood groen blauw)", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute("RGB_LED", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Font.Color = $green

# 3) "De REED magneetsensor: PC3" paragraph turns from red to green.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "De REED magneetsensor*") {
        $p.Range.Font.Color = $green
        break
    }
}
